{"js": "const replacements = [\n  [\"2026-02-26 Thursday\", \"2026-02-27 Friday\"],\n  [\"936\u00d72=\", \"965\u00d72=\"],\n  [\"683\u00d79=\", \"933\u00d76=\"],\n  [\"394\u00d76=\", \"336\u00d75=\"],\n  [\"395\u00d77=\", \"723\u00d77=\"],\n  [\"856\u00d78=\", \"574\u00d79=\"],\n  [\"881\u00d75=\", \"482\u00d75=\"],\n  [\"493\u00d72=\", \"764\u00d74=\"],\n  [\"234\u00d76=\", \"337\u00d76=\"],\n  [\"642\u00d72=\", \"707\u00d76=\"],\n  [\"420\u00d78=\", \"104\u00d75=\"],\n  [\"973\u00d79=\", \"574\u00d78=\"],\n  [\"202\u00d73=\", \"788\u00d72=\"],\n  [\"966\u00d79=\", \"360\u00d76=\"],\n  [\"599\u00d73=\", \"800\u00d79=\"],\n  [\"970\u00d75=\", \"181\u00d77=\"],\n  [\"804\u00d76=\", \"781\u00d75=\"],\n  [\"951\u00d77=\", \"459\u00d78=\"],\n  [\"922\u00d74=\", \"199\u00d75=\"],\n  [\"734\u00d77=\", \"933\u00d73=\"],\n  [\"308\u00d72=\", \"618\u00d73=\"],\n  [\"967\u00d76=\", \"121\u00d76=\"],\n  [\"523\u00d79=\", \"665\u00d75=\"],\n  [\"230\u00d76=\", \"941\u00d77=\"],\n  [\"185\u00d75=\", \"937\u00d77=\"],\n  [\"550\u00d77=\", \"279\u00d72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$replacements = @(\n    @(\"2026-02-26 Thursday\", \"2026-02-27 Friday\"),\n    @(\"936\u00d72=\", \"965\u00d72=\"),\n    @(\"683\u00d79=\", \"933\u00d76=\"),\n    @(\"394\u00d76=\", \"336\u00d75=\"),\n    @(\"395\u00d77=\", \"723\u00d77=\"),\n    @(\"856\u00d78=\", \"574\u00d79=\"),\n    @(\"881\u00d75=\", \"482\u00d75=\"),\n    @(\"493\u00d72=\", \"764\u00d74=\"),\n    @(\"234\u00d76=\", \"337\u00d76=\"),\n    @(\"642\u00d72=\", \"707\u00d76=\"),\n    @(\"420\u00d78=\", \"104\u00d75=\"),\n    @(\"973\u00d79=\", \"574\u00d78=\"),\n    @(\"202\u00d73=\", \"788\u00d72=\"),\n    @(\"966\u00d79=\", \"360\u00d76=\"),\n    @(\"599\u00d73=\", \"800\u00d79=\"),\n    @(\"970\u00d75=\", \"181\u00d77=\"),\n    @(\"804\u00d76=\", \"781\u00d75=\"),\n    @(\"951\u00d77=\", \"459\u00d78=\"),\n    @(\"922\u00d74=\", \"199\u00d75=\"),\n    @(\"734\u00d77=\", \"933\u00d73=\"),\n    @(\"308\u00d72=\", \"618\u00d73=\"),\n    @(\"967\u00d76=\", \"121\u00d76=\"),\n    @(\"523\u00d79=\", \"665\u00d75=\"),\n    @(\"230\u00d76=\", \"941\u00d77=\"),\n    @(\"185\u00d75=\", \"937\u00d77=\"),\n    @(\"550\u00d77=\", \"279\u00d72=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = $wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute([ref]$oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, [ref]$newText, $wdReplaceAll) | Out-Null\n}\n"}
